# Append two new days of commodity quotes (rows 64-65) below the existing
# data table, matching the formatting already used for the prior rows
# (date column keeps its bordered/bold/centered date-time number style).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the last existing data row onto the two new rows
# before filling in their values, so the date cells (A64/A65) keep the
# same style (bold, centered, bordered, custom date format) as A2:A63.
$ws.Range("A63:E63").Copy()
$ws.Range("A64:E65").PasteSpecial(-4122)

# Row 64 - 2022-08-01
$ws.Range("A64").Value = 44774
$ws.Range("B64").Value = 18.79000091552734
$ws.Range("C64").Value = 22.25
$ws.Range("D64").Value = 29.81999969482422
$ws.Range("E64").Value = 75.47000122070312

# Row 65 - 2022-08-02
$ws.Range("A65").Value = 44775
$ws.Range("B65").Value = 18.63500022888184
$ws.Range("C65").Value = 22.24990081787109
$ws.Range("D65").Value = 29.19000053405762
$ws.Range("E65").Value = 75.91000366210938
